# refactor(translation): added manually verified translation
#
# The English placeholder copy "Start typing here" is replaced with the
# manually verified "Start typing here..." (cell B15). Two HTML snippets
# used in the Terms & Conditions strings (B21 / D21) had their escaped
# quotes (\") cleaned up to plain quotes ("). Columns B and D are resized
# to fit their (now much longer) contents, and the active selection is
# moved to B15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Clean up the escaped quotes inside the Terms & Conditions anchor
#    markup (English + Marathi) - \" -> "
$ws.Range("B21").Value = 'By proceeding ahead you agree to the <a href="../terms-and-conditions.html" target="_blank"> Terms and Conditions</a>'
$ws.Range("D21").Value = 'पुढील कार्यवाहीसाठी तुम्हांला<a href="../terms-and-conditions.html" target="_blank"> अटीं आणि शर्तीं</a> मान्य आहे असे समजण्यात येईल'

# 2) Manually verified translation: placeholder copy gets an ellipsis.
$ws.Range("B15").Value = "Start typing here..."

# 3) Widen columns B and D so the longer copy fits.
$ws.Columns("B").ColumnWidth = 99
$ws.Columns("D").ColumnWidth = 129.75

# 4) Leave the selection on the cell that was edited.
$ws.Range("B15").Select() | Out-Null
